# Apply the header renames + drop the unused "weight by population" column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header cells (A1:D1)
$ws.Range("A1").Value = "rng_id"
$ws.Range("B1").Value = "rgn_name"
$ws.Range("C1").Value = "Year"
$ws.Range("D1").Value = "percent_households_fish"

# Remove the now-unused column E ("Weight by population to get need by island?")
$ws.Range("E1:E41").ClearContents()

# Selection moves to E1 (first empty cell after the now-cleared column)
$ws.Range("E1").Select()

$wb.Save()
